# fix: fix a PCB bug
# Seven capacitor designators (C19, C30, C45, C60, C75, C90, C105) were
# mis-assigned to the "104" (100nF) 0402-cap BOM line. Move them to the
# "10uF" 0402-cap BOM line where they actually belong.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("OECU_BASE")

# Row 15 ("104" group): remove the 7 mis-assigned designators.
# Leading "'" keeps these as literal text (preserves the existing
# quote-prefixed, wrap-text cell style instead of Excel re-deriving a
# "plain" style for the new string).
$ws.Range("C15").Value = "'C2, C3, C4, C5, C6, C7, C9, C11, C12, C13, C15, C16, C17, C18, C25, C29, C31, C32, C35, C44, C46, C47, C50, C51, C59, C61, C62, C63, C66, C74, C76, C77, C80, C83, C89, C91, C92, C93, C96, C104, C106, C107, C110, C113"

# Row 8 ("10uF" group): add the 7 designators to where they belong.
$ws.Range("C8").Value = "'C19, C30, C45, C60, C75, C90, C105, C26, C39, C43, C57, C72, C86, C99, C116"

# The "104" row's text is now shorter, so it wraps to fewer lines -
# shrink the row to fit.
$ws.Rows(15).RowHeight = 28.5

# Leave the cursor where the edit was made.
$null = $ws.Range("C12").Select()
